$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows (old "MuSCs" sending-cluster block, rows 8-10)
$ws.Rows("8:10").Delete()

# Recomputed TPM-based NATMI values for the Fgf18-Fgfr3 ligand-receptor pair
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf18"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.71557066666667
$ws.Range("H2").Value = 32.146712
$ws.Range("I2").Value = 0.9375025736567436
$ws.Range("J2").Value = 0.9375025736567436
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.752937333333333
$ws.Range("N2").Value = 11.258812
$ws.Range("O2").Value = 0.6855621274031838
$ws.Range("P2").Value = 0.6855621274031838
$ws.Range("Q2").Value = 40.21486520290489
$ws.Range("R2").Value = 361.933786826144
$ws.Range("S2").Value = 0.6427162588420772
$ws.Range("T2").Value = 0.6427162588420772
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf18"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.71557066666667
$ws.Range("H3").Value = 32.146712
$ws.Range("I3").Value = 0.9375025736567436
$ws.Range("J3").Value = 0.9375025736567436
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.296447666666667
$ws.Range("N3").Value = 3.889343
$ws.Range("O3").Value = 0.2368266084628361
$ws.Range("P3").Value = 0.2368266084628362
$ws.Range("Q3").Value = 13.89217658780178
$ws.Range("R3").Value = 125.029589290216
$ws.Range("S3").Value = 0.2220255549443068
$ws.Range("T3").Value = 0.2220255549443069
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf18"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.71557066666667
$ws.Range("H4").Value = 32.146712
$ws.Range("I4").Value = 0.9375025736567436
$ws.Range("J4").Value = 0.9375025736567436
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4248633333333334
$ws.Range("N4").Value = 1.27459
$ws.Range("O4").Value = 0.07761126413398003
$ws.Range("P4").Value = 0.07761126413398005
$ws.Range("Q4").Value = 4.552653072008889
$ws.Range("R4").Value = 40.97387764808001
$ws.Range("S4").Value = 0.0727607598703596
$ws.Range("T4").Value = 0.07276075987035961
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Fgf18"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.71434
$ws.Range("H5").Value = 2.14302
$ws.Range("I5").Value = 0.06249742634325634
$ws.Range("J5").Value = 0.06249742634325634
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.752937333333333
$ws.Range("N5").Value = 11.258812
$ws.Range("O5").Value = 0.6855621274031838
$ws.Range("P5").Value = 0.6855621274031838
$ws.Range("Q5").Value = 2.680873254693333
$ws.Range("R5").Value = 24.12785929224
$ws.Range("S5").Value = 0.0428458685611066
$ws.Range("T5").Value = 0.0428458685611066
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Fgf18"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.71434
$ws.Range("H6").Value = 2.14302
$ws.Range("I6").Value = 0.06249742634325634
$ws.Range("J6").Value = 0.06249742634325634
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.296447666666667
$ws.Range("N6").Value = 3.889343
$ws.Range("O6").Value = 0.2368266084628361
$ws.Range("P6").Value = 0.2368266084628362
$ws.Range("Q6").Value = 0.9261044262066667
$ws.Range("R6").Value = 8.33493983586
$ws.Range("S6").Value = 0.01480105351852931
$ws.Range("T6").Value = 0.01480105351852931
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Fgf18"
$ws.Range("C7").Value = "Fgfr3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.71434
$ws.Range("H7").Value = 2.14302
$ws.Range("I7").Value = 0.06249742634325634
$ws.Range("J7").Value = 0.06249742634325634
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4248633333333334
$ws.Range("N7").Value = 1.27459
$ws.Range("O7").Value = 0.07761126413398003
$ws.Range("P7").Value = 0.07761126413398005
$ws.Range("Q7").Value = 0.3034968735333333
$ws.Range("R7").Value = 2.7314718618
$ws.Range("S7").Value = 0.00485050426362043
$ws.Range("T7").Value = 0.00485050426362043
